$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44208
$ws.Range("J2").Value = 85
$ws.Range("K2").Value = 3700
$ws.Range("M2").Value = 3824
$ws.Range("P2").Value = 1912

# Row 3
$ws.Range("D3").Value = 44664

# Row 4
$ws.Range("D4").Value = 44662
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8500
$ws.Range("M4").Value = 8250
$ws.Range("P4").Value = 229

# Row 5
$ws.Range("D5").Value = 44225
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 3400
$ws.Range("L5").Value = 3700
$ws.Range("M5").Value = 3550
$ws.Range("P5").Value = 1775

# Row 6
$ws.Range("D6").Value = 44161
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 2800
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 2900
$ws.Range("P6").Value = 1450

# Row 7 - unchanged

# Row 8
$ws.Range("D8").Value = 44166
$ws.Range("J8").Value = 70
$ws.Range("K8").Value = 3500
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3679
$ws.Range("N8").Value = "$/paquete 36 unidades"
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 102
$ws.Range("Q8").Value = 36

# Row 9
$ws.Range("D9").Value = 44215
$ws.Range("J9").Value = 140
$ws.Range("M9").Value = 3768
$ws.Range("P9").Value = 1884

# Row 10
$ws.Range("D10").Value = 44160
$ws.Range("J10").Value = 43
$ws.Range("K10").Value = 3500
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = 3709
$ws.Range("P10").Value = 103

# Row 11
$ws.Range("D11").Value = 44209
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 3500
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 3767
$ws.Range("P11").Value = 1884

# Row 12
$ws.Range("D12").Value = 44210
$ws.Range("J12").Value = 105
$ws.Range("M12").Value = 3714
$ws.Range("N12").Value = "$/paquete 2 kilos"
$ws.Range("O12").Value = "Provincia de Diguillín"
$ws.Range("P12").Value = 1857
$ws.Range("Q12").Value = 2
